$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 443, shifting existing rows 443:457 down to 444:458
$ws.Rows("443:443").Insert()

# Populate the newly inserted row 443 with the new data record
$ws.Range("A443").Value = 4
$ws.Range("B443").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C443").Value = "Los Lagos"
$ws.Range("D443").Value = 45075
$ws.Range("E443").Value = 10
$ws.Range("F443").Value = 100112037
$ws.Range("G443").Value = "Cebollín"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 70
$ws.Range("K443").Value = 7000
$ws.Range("L443").Value = 7000
$ws.Range("M443").Value = 7000
$ws.Range("N443").Value = "$/paquete 36 unidades"
$ws.Range("O443").Value = "Región Metropolitana"
$ws.Range("P443").Value = 194
$ws.Range("Q443").Value = 36
$ws.Range("R443").Value = "Hortaliza"
